$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 16426
$ws.Range("J32").Value = 19002
$ws.Range("L32").Value = 19002
$ws.Range("N32").Value = -19654
$ws.Range("H98").Value = 224018.73
$ws.Range("I98").Value = 937.375
$ws.Range("K98").Value = 937.375
$ws.Range("M98").Value = 560.625
$ws.Range("H116").Value = 8538
$ws.Range("I116").Value = 7198.2856
$ws.Range("K116").Value = 7198.2856
$ws.Range("M116").Value = -3756.2856
$ws.Range("H122").Value = 224018.73
$ws.Range("I122").Value = 937.375
$ws.Range("K122").Value = 2812.125
$ws.Range("M122").Value = -362.125
$ws.Range("H132").Value = 3734.182
$ws.Range("I132").Value = 3734.182
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11202.546
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 3636.55
$ws.Range("J138").Value = 4019.9539
$ws.Range("L138").Value = 12059.8617
$ws.Range("N138").Value = -22339.8617

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 287
$ws.Range("I5").Value = 222.1
$ws.Range("J5").Value = 449.25
$ws.Range("K5").Value = 222.1
$ws.Range("L5").Value = 449.25
$ws.Range("M5").Value = -110.1
$ws.Range("N5").Value = -673.25
$ws.Range("H32").Value = 5567.4463
$ws.Range("I32").Value = 5868.1763
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 5868.1763
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = -5581.1763
$ws.Range("N32").Value = -3074
$ws.Range("H45").Value = 4049.2856
$ws.Range("I45").Value = 3629.111
$ws.Range("K45").Value = 3629.111
$ws.Range("M45").Value = -3252.111
$ws.Range("H97").Value = 1240.5769
$ws.Range("I97").Value = 1034.3636
$ws.Range("K97").Value = 1034.3636
$ws.Range("M97").Value = -538.3635999999999
$ws.Range("H122").Value = 3692.7083
$ws.Range("I122").Value = 3477.4167
$ws.Range("J122").Value = 3908
$ws.Range("K122").Value = 10432.2501
$ws.Range("L122").Value = 11724
$ws.Range("M122").Value = -7982.250100000001
$ws.Range("N122").Value = -16624
$ws.Range("H132").Value = 3709.1365
$ws.Range("J132").Value = 10506.5
$ws.Range("L132").Value = 31519.5
$ws.Range("N132").Value = -36579.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 287
$ws.Range("I4").Value = 222.1
$ws.Range("J4").Value = 449.25
$ws.Range("K4").Value = 222.1
$ws.Range("L4").Value = 449.25
$ws.Range("M4").Value = -107.1
$ws.Range("N4").Value = -679.25
$ws.Range("H94").Value = 976.8723
$ws.Range("I94").Value = 1049.475
$ws.Range("K94").Value = 1049.475
$ws.Range("M94").Value = -598.4749999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 349
$ws.Range("I7").Value = 138.33333
$ws.Range("J7").Value = 529.5714
$ws.Range("K7").Value = 138.33333
$ws.Range("L7").Value = 529.5714
$ws.Range("M7").Value = -25.33332999999999
$ws.Range("N7").Value = -755.5714
$ws.Range("H22").Value = 4680.1816
$ws.Range("I22").Value = 1126.8
$ws.Range("K22").Value = 1126.8
$ws.Range("M22").Value = -776.8
$ws.Range("H31").Value = 31057.621
$ws.Range("I31").Value = 1967.6923
$ws.Range("J31").Value = 99815.63
$ws.Range("K31").Value = 1967.6923
$ws.Range("L31").Value = 99815.63
$ws.Range("M31").Value = -1672.6923
$ws.Range("N31").Value = -100405.63
$ws.Range("H34").Value = 31057.621
$ws.Range("I34").Value = 1967.6923
$ws.Range("J34").Value = 99815.63
$ws.Range("K34").Value = 1967.6923
$ws.Range("L34").Value = 99815.63
$ws.Range("M34").Value = -1765.6923
$ws.Range("N34").Value = -100219.63
$ws.Range("H107").Value = 2283.238
$ws.Range("I107").Value = 1713.0834
$ws.Range("K107").Value = 1713.0834
$ws.Range("M107").Value = 206.9166
$ws.Range("H134").Value = 3314.7058
$ws.Range("I134").Value = 2226.9167
$ws.Range("K134").Value = 6680.750100000001
$ws.Range("M134").Value = -4145.750100000001
$ws.Range("H135").Value = 55905.91
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 55905.91
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -66045.91

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 71719.86
$ws.Range("I7").Value = 422.25
$ws.Range("K7").Value = 1266.75
$ws.Range("M7").Value = -1154.75
$ws.Range("H81").Value = 3862.3333
$ws.Range("I81").Value = 219.66667
$ws.Range("J81").Value = 7505
$ws.Range("K81").Value = 659.00001
$ws.Range("L81").Value = 22515
$ws.Range("M81").Value = 463.99999
$ws.Range("N81").Value = -24761
$ws.Range("H84").Value = 3862.3333
$ws.Range("I84").Value = 219.66667
$ws.Range("J84").Value = 7505
$ws.Range("K84").Value = 1977.00003
$ws.Range("L84").Value = 67545
$ws.Range("M84").Value = 3638.99997
$ws.Range("N84").Value = -78777
$ws.Range("H107").Value = 1489131.5
$ws.Range("I107").Value = 921.25
$ws.Range("J107").Value = 2404953.2
$ws.Range("K107").Value = 2763.75
$ws.Range("L107").Value = 7214859.600000001
$ws.Range("M107").Value = -843.75
$ws.Range("N107").Value = -7218699.600000001
$ws.Range("H109").Value = 2775.4443
$ws.Range("I109").Value = 1517.5
$ws.Range("K109").Value = 4552.5
$ws.Range("M109").Value = -3512.5
$ws.Range("H114").Value = 2615.3333
$ws.Range("I114").Value = 553.6667
$ws.Range("J114").Value = 4677
$ws.Range("K114").Value = 1661.0001
$ws.Range("L114").Value = 14031
$ws.Range("M114").Value = 1592.9999
$ws.Range("N114").Value = -20539
$ws.Range("H139").Value = 2203.2258
$ws.Range("I139").Value = 1293.375
$ws.Range("J139").Value = 5322.7144
$ws.Range("K139").Value = 3880.125
$ws.Range("L139").Value = 15968.1432
$ws.Range("M139").Value = 1259.875
$ws.Range("N139").Value = -26248.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3424.9
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H80").Value = 363827.72
$ws.Range("I80").Value = 835598.5
$ws.Range("J80").Value = 9999.625
$ws.Range("K80").Value = 835598.5
$ws.Range("L80").Value = 9999.625
$ws.Range("M80").Value = -834600.5
$ws.Range("N80").Value = -11995.625
$ws.Range("H83").Value = 363827.72
$ws.Range("I83").Value = 835598.5
$ws.Range("J83").Value = 9999.625
$ws.Range("K83").Value = 4177992.5
$ws.Range("L83").Value = 49998.125
$ws.Range("M83").Value = -4173000.5
$ws.Range("N83").Value = -59982.125
$ws.Range("H99").Value = 25340.625
$ws.Range("I99").Value = 15537.5
$ws.Range("K99").Value = 15537.5
$ws.Range("M99").Value = -13291.5
$ws.Range("H102").Value = 2856.087
$ws.Range("I102").Value = 1275.1177
$ws.Range("J102").Value = 7335.5
$ws.Range("K102").Value = 1275.1177
$ws.Range("L102").Value = 7335.5
$ws.Range("M102").Value = 346.8823
$ws.Range("N102").Value = -10579.5
$ws.Range("H113").Value = 3582.762
$ws.Range("I113").Value = 2868.2222
$ws.Range("K113").Value = 2868.2222
$ws.Range("M113").Value = -698.2222000000002
$ws.Range("H122").Value = 5051.636
$ws.Range("I122").Value = 4926.923
$ws.Range("K122").Value = 14780.769
$ws.Range("M122").Value = -12330.769
$ws.Range("H132").Value = 5268.372
$ws.Range("I132").Value = 3628.9333
$ws.Range("J132").Value = 9051.691999999999
$ws.Range("K132").Value = 10886.7999
$ws.Range("L132").Value = 27155.076
$ws.Range("M132").Value = -8356.7999
$ws.Range("N132").Value = -32215.076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6200.2
$ws.Range("I7").Value = 4343.727
$ws.Range("J7").Value = 11305.5
$ws.Range("K7").Value = 4343.727
$ws.Range("L7").Value = 11305.5
$ws.Range("M7").Value = -4231.727
$ws.Range("N7").Value = -11529.5
$ws.Range("H40").Value = 6206.483
$ws.Range("I40").Value = 5557.1177
$ws.Range("K40").Value = 5557.1177
$ws.Range("M40").Value = -5421.1177
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H126").Value = 6200.2
$ws.Range("I126").Value = 4343.727
$ws.Range("J126").Value = 11305.5
$ws.Range("K126").Value = 13031.181
$ws.Range("L126").Value = 33916.5
$ws.Range("M126").Value = -10561.181
$ws.Range("N126").Value = -38856.5
$ws.Range("H136").Value = 4455.3887
$ws.Range("I136").Value = 3262
$ws.Range("K136").Value = 9786
$ws.Range("M136").Value = -7236

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 60000
$ws.Range("J46").Value = 60000
$ws.Range("L46").Value = 60000
$ws.Range("N46").Value = -60462
$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 180000
$ws.Range("N134").Value = -185070
